# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = 0.6545652718822623
    "C2" = 1.626987699542094
    "D2" = 0.1496068669990043
    "E2" = 0.5333859586016987
    "G2" = 2.964545797025059

    "B3" = 0.04172184405617529
    "C3" = 0.3048912486333797
    "D3" = 189.6080260415259
    "E3" = 13.86384647080068
    "G3" = 203.8184856050161

    "B4" = 1.445647641019636
    "C4" = 0.3048912486333797
    "D4" = 0.7210945179870265
    "E4" = 0.5333859586016987
    "G4" = 3.005019366241741

    "B5" = 3.272327238179451
    "C5" = 1.626987699542094
    "D5" = 18.71679738969934
    "E5" = 0.5333859586016987
    "G5" = 24.14949828602258

    "B6" = 3.272327238179451
    "C6" = 1.626987699542094
    "D6" = 0.7210945179870265
    "E6" = 0.5333859586016987
    "G6" = 6.15379541431027
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
